# "added popupadress as subview"
# Updates the DTT-Test-Hour-Log worksheet:
#  - Rewrites the "Amount of hours" values in column C that were stored as the
#    text strings "1.50"/"0.50" so they instead reuse the already-present
#    "1.5"/"0.5" text values (dedupes the shared-string table).
#  - Fills in the previously-blank row 24 with a new log entry.
#  - Updates the sheet's selection to the cell the author was last working on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as TEXT (not as a number) while leaving
# the cell's existing style/number-format index untouched. We do this by
# temporarily switching the cell to a text format, assigning the value, and
# then restoring the original (built-in "0") number format used by this
# column - this forces Excel to store the cell as a shared string (t="s")
# instead of coercing the numeric-looking text into a real number.
function Set-TextValue {
    param($address, $value)
    $rng = $ws.Range($address)
    $rng.NumberFormat = "@"
    $rng.Value2 = $value
    $rng.NumberFormat = "0"
}

# Row 24 was completely empty; populate it with the new log entry first so
# the two brand-new shared strings it introduces are appended to the shared
# string table ahead of the new "1.5" value below (matches the order in
# which the author's edit introduced them).
$ws.Range("A24").Value2 = "Seperated all methods into extentions"
$ws.Range("B24").Value2 = 43103
Set-TextValue "C24" "0.5"
$ws.Range("D24").Value2 = "I seperated all the methodsd into the correct extentions. "

# Column C ("Amount of hours") entries that were textual "1.50"/"0.50" are
# normalized to "1.5"/"0.5" (same text cells used elsewhere in the sheet).
Set-TextValue "C14" "0.5"
Set-TextValue "C15" "1.5"
Set-TextValue "C16" "1.5"
Set-TextValue "C17" "1.5"
Set-TextValue "C18" "0.5"
Set-TextValue "C23" "0.5"

# Reflect the author's last viewport/selection in the sheet view state.
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("C18").Select()
